$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two additional example data cells (scoreAwarded / qualificationAwarded)
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = "A"

# Move the selection / viewport the way the author left it after editing
$ws.Range("M8").Select()
